# Adds a "Persoon" column (D) to the opdrachten schema, assigning each
# opdracht to Jeroen or Laura, and updates the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("D3").Value = "Persoon"

# Assign a person to each "opdracht" row
$ws.Range("D4").Value  = "Jeroen"
$ws.Range("D5").Value  = "Jeroen"
$ws.Range("D6").Value  = "Laura"
$ws.Range("D7").Value  = "Laura"
$ws.Range("D8").Value  = "Laura"
$ws.Range("D9").Value  = "Jeroen"
$ws.Range("D10").Value = "Laura"
$ws.Range("D11").Value = "Jeroen"

# Match the updated view state: scrolled one column right, A11 selected
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A11").Select()
